$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Value edits to existing cells (F column purchase prices + L column extra
#    costs). H column formulas recalc automatically.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = 410000
$ws.Range("F3").Value = 513400
$ws.Range("F5").Value = 20000
$ws.Range("F7").Value = 1820000

$ws.Range("L3").Value = 300000
$ws.Range("L4").Value = 70000
$ws.Range("L5").Value = 10000
$ws.Range("L6").Value = 60000

# ---------------------------------------------------------------------------
# 2) Move the supplier-link text currently sitting in N3 over to M3, clearing
#    N3 out (it keeps a plain, non-wrapping style).
# ---------------------------------------------------------------------------
$n3Text = $ws.Range("N3").Value
$ws.Range("N3").Value = ""
$ws.Range("N3").WrapText = $false

# ---------------------------------------------------------------------------
# 3) Apply the "note" style (font/border copied from column B, wrap text
#    turned off) to the whole M2:M8 helper column, then fill in the two rows
#    that carry actual notes.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("M2").WrapText = $false

$ws.Range("B3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").WrapText = $false
$ws.Range("M3").Value = $n3Text

$ws.Range("B4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").WrapText = $false

$ws.Range("B5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").WrapText = $false
$ws.Range("M5").Value = "95 450 02 00 Shukrillo"

$ws.Range("B6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").WrapText = $false

$ws.Range("B7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").WrapText = $false

$ws.Range("B8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").WrapText = $false

# Widen column M to fit the longer note text.
$ws.Columns.Item(13).ColumnWidth = 42.42578125

# ---------------------------------------------------------------------------
# 4) New hyperlink in N5 pointing at the olx.uz listing.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("N5"), "https://www.olx.uz/d/obyavlenie/trelleborg-uplotnitel-samokleyuschiysya-ID3sPb4.html?reason=extended_search_extended_category")

# ---------------------------------------------------------------------------
# 5) Grand-total row under the table.
# ---------------------------------------------------------------------------
$ws.Range("H8").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Formula = "=SUM(H2:H9)"

# ---------------------------------------------------------------------------
# 6) Selection cursor, matching where the editor left off.
# ---------------------------------------------------------------------------
$ws.Range("F9").Select()
